$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated mat_rank (G), index (C) and race-string (H) values per row.
# Race string values map to shared strings: White, Asian, Black or African American, Hispanic
$raceWhite = "White"
$raceAsian = "Asian"
$raceBlack = "Black or African American"
$raceHispanic = "Hispanic"

# Row 2
$ws.Range("G2").Value = 13.4562512144795

# Row 3
$ws.Range("G3").Value = 13.03225951465968

# Row 4
$ws.Range("C4").Value = 21
$ws.Range("G4").Value = 8.49956431831203
$ws.Range("H4").Value = $raceBlack

# Row 5
$ws.Range("G5").Value = 8.20087022370102

# Row 6
$ws.Range("C6").Value = 22
$ws.Range("G6").Value = 8.064168822103696
$ws.Range("H6").Value = $raceAsian

# Row 7
$ws.Range("G7").Value = 5.333931338090698

# Row 8
$ws.Range("G8").Value = 5.316202313826643

# Row 9
$ws.Range("G9").Value = 5.294121455295787

# Row 10
$ws.Range("C10").Value = 34
$ws.Range("G10").Value = 4.302967855272668

# Row 11
$ws.Range("C11").Value = 35
$ws.Range("G11").Value = 4.23075704731449

# Row 12
$ws.Range("G12").Value = 2.016984074606204

# Row 13
$ws.Range("G13").Value = 1.054877676087834

# Row 14
$ws.Range("G14").Value = 14.32124806351207

# Row 15
$ws.Range("G15").Value = 13.25929927578149

# Row 16
$ws.Range("G16").Value = 8.119568600285705

# Row 17
$ws.Range("G17").Value = 7.411316600731239

# Row 18
$ws.Range("G18").Value = 6.387209556654361

# Row 19
$ws.Range("G19").Value = 6.159915716578424

# Row 20
$ws.Range("C20").Value = 32
$ws.Range("G20").Value = 5.496872041548905
$ws.Range("H20").Value = $raceBlack

# Row 21
$ws.Range("C21").Value = 30
$ws.Range("G21").Value = 5.269135301367183

# Row 22
$ws.Range("C22").Value = 33
$ws.Range("G22").Value = 5.262741384947466
$ws.Range("H22").Value = $raceWhite

# Row 23
$ws.Range("G23").Value = 3.214622224054206

# Row 24
$ws.Range("G24").Value = 1.200783564827204

# Row 25
$ws.Range("G25").Value = 0.09961791273931464
